# Update values on Sheet1 as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.228
$ws.Range("A9").Value = -20.912
$ws.Range("B11").Value = 6.927
$ws.Range("A18").Value = -21.825
$ws.Range("A20").Value = -21.738
$ws.Range("D21").Value = -7.675999999999999
